# Scheduled-runner refresh of market-price-derived columns
# (currentAveragePrice*, Leve*Price*, Leve*Profit*) across the per-job
# leve-profit tables. Source game-item IDs/dates are unchanged; only the
# price-derived columns (H:N) are refreshed with newly pulled data.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 9549.15
$ws.Range("I28").Value = 1586.6666
$ws.Range("J28").Value = 33436.6
$ws.Range("K28").Value = 1586.6666
$ws.Range("L28").Value = 33436.6
$ws.Range("M28").Value = -1101.6666
$ws.Range("N28").Value = -34406.6
$ws.Range("H62").Value = 6927.2856
$ws.Range("I62").Value = 2322.875
$ws.Range("K62").Value = 2322.875
$ws.Range("M62").Value = -1698.875
$ws.Range("H65").Value = 6927.2856
$ws.Range("I65").Value = 2322.875
$ws.Range("K65").Value = 11614.375
$ws.Range("M65").Value = -8494.375
$ws.Range("H98").Value = 988.7037
$ws.Range("I98").Value = 847.96
$ws.Range("J98").Value = 2748
$ws.Range("K98").Value = 847.96
$ws.Range("L98").Value = 2748
$ws.Range("M98").Value = 650.04
$ws.Range("N98").Value = -5744
$ws.Range("H107").Value = 1380.1305
$ws.Range("I107").Value = 1757
$ws.Range("J107").Value = 793.8889
$ws.Range("K107").Value = 1757
$ws.Range("L107").Value = 793.8889
$ws.Range("M107").Value = 163
$ws.Range("N107").Value = -4633.8889
$ws.Range("H122").Value = 988.7037
$ws.Range("I122").Value = 847.96
$ws.Range("J122").Value = 2748
$ws.Range("K122").Value = 2543.88
$ws.Range("L122").Value = 8244
$ws.Range("M122").Value = -93.88000000000011
$ws.Range("N122").Value = -13144
$ws.Range("H129").Value = 1417.289
$ws.Range("I129").Value = 446.53845
$ws.Range("J129").Value = 1811.6562
$ws.Range("K129").Value = 1339.61535
$ws.Range("L129").Value = 5434.9686
$ws.Range("M129").Value = 3660.38465
$ws.Range("N129").Value = -15434.9686
$ws.Range("H132").Value = 3937.6086
$ws.Range("I132").Value = 4162.4736
$ws.Range("J132").Value = 2869.5
$ws.Range("K132").Value = 12487.4208
$ws.Range("L132").Value = 8608.5
$ws.Range("M132").Value = -9957.4208
$ws.Range("N132").Value = -13668.5
$ws.Range("H137").Value = 1899.5927
$ws.Range("I137").Value = 1991.5454
$ws.Range("J137").Value = 1836.375
$ws.Range("K137").Value = 5974.6362
$ws.Range("L137").Value = 5509.125
$ws.Range("M137").Value = -3424.6362
$ws.Range("N137").Value = -10609.125

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2752
$ws.Range("I45").Value = 2800
$ws.Range("J45").Value = 2704
$ws.Range("K45").Value = 2800
$ws.Range("L45").Value = 2704
$ws.Range("M45").Value = -2423
$ws.Range("N45").Value = -3458
$ws.Range("H61").Value = 8337294
$ws.Range("I61").Value = 16668983
$ws.Range("J61").Value = 5605.6
$ws.Range("K61").Value = 16668983
$ws.Range("L61").Value = 5605.6
$ws.Range("M61").Value = -16668771
$ws.Range("N61").Value = -6029.6
$ws.Range("H74").Value = 10418972
$ws.Range("I74").Value = 1363.6061
$ws.Range("J74").Value = 33337712
$ws.Range("K74").Value = 1363.6061
$ws.Range("L74").Value = 33337712
$ws.Range("M74").Value = -489.6061
$ws.Range("N74").Value = -33339460
$ws.Range("H77").Value = 10418972
$ws.Range("I77").Value = 1363.6061
$ws.Range("J77").Value = 33337712
$ws.Range("K77").Value = 6818.0305
$ws.Range("L77").Value = 166688560
$ws.Range("M77").Value = -2450.0305
$ws.Range("N77").Value = -166697296
$ws.Range("H136").Value = 8337294
$ws.Range("I136").Value = 16668983
$ws.Range("J136").Value = 5605.6
$ws.Range("K136").Value = 50006949
$ws.Range("L136").Value = 16816.8
$ws.Range("M136").Value = -50004399
$ws.Range("N136").Value = -21916.8

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2625.0278
$ws.Range("I134").Value = 2453.1562
$ws.Range("J134").Value = 4000
$ws.Range("K134").Value = 7359.4686
$ws.Range("L134").Value = 12000
$ws.Range("M134").Value = -4824.4686
$ws.Range("N134").Value = -17070

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4802.2856
$ws.Range("I31").Value = 1301.1852
$ws.Range("J31").Value = 7000.6514
$ws.Range("K31").Value = 1301.1852
$ws.Range("L31").Value = 7000.6514
$ws.Range("M31").Value = -1006.1852
$ws.Range("N31").Value = -7590.6514
$ws.Range("H34").Value = 4802.2856
$ws.Range("I34").Value = 1301.1852
$ws.Range("J34").Value = 7000.6514
$ws.Range("K34").Value = 1301.1852
$ws.Range("L34").Value = 7000.6514
$ws.Range("M34").Value = -1099.1852
$ws.Range("N34").Value = -7404.6514
$ws.Range("H58").Value = 2972.3684
$ws.Range("I58").Value = 2811.5
$ws.Range("J58").Value = 3422.8
$ws.Range("K58").Value = 2811.5
$ws.Range("L58").Value = 3422.8
$ws.Range("M58").Value = -2608.5
$ws.Range("N58").Value = -3828.8
$ws.Range("H134").Value = 8338896
$ws.Range("I134").Value = 11910899
$ws.Range("J134").Value = 4222.1113
$ws.Range("K134").Value = 35732697
$ws.Range("L134").Value = 12666.3339
$ws.Range("M134").Value = -35730162
$ws.Range("N134").Value = -17736.3339
$ws.Range("H136").Value = 2972.3684
$ws.Range("I136").Value = 2811.5
$ws.Range("J136").Value = 3422.8
$ws.Range("K136").Value = 8434.5
$ws.Range("L136").Value = 10268.4
$ws.Range("M136").Value = -5884.5
$ws.Range("N136").Value = -15368.4

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H98").Value = 211466.67
$ws.Range("J98").Value = 271685.72
$ws.Range("L98").Value = 815057.1599999999
$ws.Range("N98").Value = -818053.1599999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 11906
$ws.Range("I70").Value = 15298.223
$ws.Range("J70").Value = 5800
$ws.Range("K70").Value = 15298.223
$ws.Range("L70").Value = 5800
$ws.Range("M70").Value = -15028.223
$ws.Range("N70").Value = -6340
$ws.Range("H73").Value = 11906
$ws.Range("I73").Value = 15298.223
$ws.Range("J73").Value = 5800
$ws.Range("K73").Value = 15298.223
$ws.Range("L73").Value = 5800
$ws.Range("M73").Value = -14362.223
$ws.Range("N73").Value = -7672
$ws.Range("H80").Value = 785002.4399999999
$ws.Range("I80").Value = 3001335
$ws.Range("J80").Value = 46224.89
$ws.Range("K80").Value = 3001335
$ws.Range("L80").Value = 46224.89
$ws.Range("M80").Value = -3000337
$ws.Range("N80").Value = -48220.89
$ws.Range("H83").Value = 785002.4399999999
$ws.Range("I83").Value = 3001335
$ws.Range("J83").Value = 46224.89
$ws.Range("K83").Value = 15006675
$ws.Range("L83").Value = 231124.45
$ws.Range("M83").Value = -15001683
$ws.Range("N83").Value = -241108.45
$ws.Range("H123").Value = 8991.066000000001
$ws.Range("J123").Value = 8991.066000000001
$ws.Range("L123").Value = 8991.066000000001
$ws.Range("N123").Value = -13891.066

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 3321.625
$ws.Range("I61").Value = 2467.818
$ws.Range("J61").Value = 5200
$ws.Range("K61").Value = 2467.818
$ws.Range("L61").Value = 5200
$ws.Range("M61").Value = -2265.818
$ws.Range("N61").Value = -5604
$ws.Range("H113").Value = 3321.625
$ws.Range("I113").Value = 2467.818
$ws.Range("J113").Value = 5200
$ws.Range("K113").Value = 2467.818
$ws.Range("L113").Value = 5200
$ws.Range("M113").Value = -297.8180000000002
$ws.Range("N113").Value = -9540
$ws.Range("H132").Value = 3234.121
$ws.Range("I132").Value = 2775.2942
$ws.Range("J132").Value = 3721.625
$ws.Range("K132").Value = 8325.882599999999
$ws.Range("L132").Value = 11164.875
$ws.Range("M132").Value = -5795.882599999999
$ws.Range("N132").Value = -16224.875

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 103246.08
$ws.Range("I62").Value = 127587.375
$ws.Range("J62").Value = 64300
$ws.Range("K62").Value = 127587.375
$ws.Range("L62").Value = 64300
$ws.Range("M62").Value = -126963.375
$ws.Range("N62").Value = -65548
$ws.Range("H65").Value = 103246.08
$ws.Range("I65").Value = 127587.375
$ws.Range("J65").Value = 64300
$ws.Range("K65").Value = 637936.875
$ws.Range("L65").Value = 321500
$ws.Range("M65").Value = -634816.875
$ws.Range("N65").Value = -327740
